$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare style for new rows 17-19 column A (copy formatting from an existing styled cell)
$ws.Range("A10").Copy($ws.Range("A17"))
$ws.Range("A10").Copy($ws.Range("A18"))
$ws.Range("A10").Copy($ws.Range("A19"))

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.007487963729129
$ws.Range("D10").Value = 0.9397300125747868
$ws.Range("E10").Value = 1.009179549290269
$ws.Range("F10").Value = 1.007487963729129
$ws.Range("G10").Value = 0.9652064309417071
$ws.Range("H10").Value = 1.031576304979364
$ws.Range("I10").Value = 1.007988454637791
$ws.Range("J10").Value = 0.9397300125747868
$ws.Range("K10").Value = 0.9744547809325279
$ws.Range("L10").Value = 0.9909713723308284
$ws.Range("M10").Value = 0.9935281193588413

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9979308762063809
$ws.Range("D11").Value = 0.9661648987374546
$ws.Range("E11").Value = 1.005075894652692
$ws.Range("F11").Value = 0.9979308762063809
$ws.Range("G11").Value = 0.9785583953420379
$ws.Range("H11").Value = 1.021637462913373
$ws.Range("I11").Value = 1.002417402898348
$ws.Range("J11").Value = 0.9661648987374546
$ws.Range("K11").Value = 0.9856203966950732
$ws.Range("L11").Value = 0.991775636450727
$ws.Range("M11").Value = 0.9952974884583811

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9978254906456129
$ws.Range("D12").Value = 0.9667318551334873
$ws.Range("E12").Value = 1.004936322300331
$ws.Range("F12").Value = 0.9978254906456129
$ws.Range("G12").Value = 0.9788606202377387
$ws.Range("H12").Value = 1.02119392308064
$ws.Range("I12").Value = 1.002292522245682
$ws.Range("J12").Value = 0.9667318551334873
$ws.Range("K12").Value = 0.985834088716909
$ws.Range("L12").Value = 0.9918297896812609
$ws.Range("M12").Value = 0.9953067889405821

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9978885232742185
$ws.Range("D13").Value = 0.9662763215043029
$ws.Range("E13").Value = 1.005053398882734
$ws.Range("F13").Value = 0.9978885232742185
$ws.Range("G13").Value = 0.9786078552373842
$ws.Range("H13").Value = 1.021549924573545
$ws.Range("I13").Value = 1.002392767978365
$ws.Range("J13").Value = 0.9662763215043029
$ws.Range("K13").Value = 0.9856648601935183
$ws.Range("L13").Value = 0.9917766917338684
$ws.Range("M13").Value = 0.9952947985750917

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.021475999999999
$ws.Range("D14").Value = 0.890344
$ws.Range("E14").Value = 1.022724
$ws.Range("F14").Value = 1.021475999999999
$ws.Range("G14").Value = 0.9393400000000003
$ws.Range("H14").Value = 1.072339999999998
$ws.Range("I14").Value = 1.021687999999998
$ws.Range("J14").Value = 0.890344
$ws.Range("K14").Value = 0.9565339999999998
$ws.Range("L14").Value = 0.9890049999999996
$ws.Range("M14").Value = 0.9946519999999993

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.05
$ws.Range("D15").Value = 0.78
$ws.Range("E15").Value = 1.05
$ws.Range("F15").Value = 1.05
$ws.Range("G15").Value = 0.88
$ws.Range("H15").Value = 1.15
$ws.Range("I15").Value = 1.05
$ws.Range("J15").Value = 0.78
$ws.Range("K15").Value = 0.915
$ws.Range("L15").Value = 0.9824999999999999
$ws.Range("M15").Value = 0.9933333333333333

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.0268011085824
$ws.Range("D16").Value = 0.8697213132800017
$ws.Range("E16").Value = 1.027072780800001
$ws.Range("F16").Value = 1.0268011085824
$ws.Range("G16").Value = 0.928115720806401
$ws.Range("H16").Value = 1.084378283622395
$ws.Range("I16").Value = 1.026701848985601
$ws.Range("J16").Value = 0.8697213132800017
$ws.Range("K16").Value = 0.9483970470400014
$ws.Range("L16").Value = 0.9875990778112006
$ws.Range("M16").Value = 0.9937985093461332

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.994581720266642
$ws.Range("D17").Value = 0.9951276124217265
$ws.Range("E17").Value = 0.9948179213860001
$ws.Range("F17").Value = 0.994581720266642
$ws.Range("G17").Value = 0.9947283782026042
$ws.Range("H17").Value = 0.9952215673587602
$ws.Range("I17").Value = 0.9947889018472132
$ws.Range("J17").Value = 0.9951276124217265
$ws.Range("K17").Value = 0.9949727669038633
$ws.Range("L17").Value = 0.9947772435852527
$ws.Range("M17").Value = 0.9948776835804911

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9931209534802544
$ws.Range("D18").Value = 1.002662861815081
$ws.Range("E18").Value = 0.9924736702166709
$ws.Range("F18").Value = 0.9931209534802544
$ws.Range("G18").Value = 0.9984126915535669
$ws.Range("H18").Value = 0.9913181358439231
$ws.Range("I18").Value = 0.9913486630203208
$ws.Range("J18").Value = 1.002662861815081
$ws.Range("K18").Value = 0.997568266015876
$ws.Range("L18").Value = 0.9953446097480653
$ws.Range("M18").Value = 0.9948894959883029

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9898988303559987
$ws.Range("D19").Value = 1.018046544847298
$ws.Range("E19").Value = 0.9890910863222876
$ws.Range("F19").Value = 0.9898988303559987
$ws.Range("G19").Value = 1.008742264628955
$ws.Range("H19").Value = 0.9783738948154228
$ws.Range("I19").Value = 0.9886187749457014
$ws.Range("J19").Value = 1.018046544847298
$ws.Range("K19").Value = 1.003568815584793
$ws.Range("L19").Value = 0.9967338229703959
$ws.Range("M19").Value = 0.9954618993192774
